# Scrape-run update for línea 141 (31/12/2025 10:42): adds newly scraped
# arrival rows to each of the three sheets and refreshes the "Última
# actualización" / "Total filas" banner cells (A2 / A3) on every sheet.

$wb = $excel.ActiveWorkbook

$updated = "Última actualización: 31/12/2025 10:42:51"

# ---------------------------------------------------------------------
# Sheet "LP1912": columns A(-) B=Hora_Scrap C=Hora_Llegada D=Linea
#                 E=Minutos F=Parada G=Fecha
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = $updated
$ws1.Range("A3").Value = "Total filas: 851"

$rows1 = @(
    @(833, "10:42:40", "10:44", "14_ABASTO",          2, "LP1912", "31/12/2025"),
    @(834, "10:42:40", "10:49", "16_SANTA ANA",        7, "LP1912", "31/12/2025"),
    @(835, "10:42:40", "10:51", "15_ABASTO",           9, "LP1912", "31/12/2025"),
    @(836, "10:42:40", "10:54", "10_OLMOS",           12, "LP1912", "31/12/2025"),
    @(837, "10:42:40", "10:56", "27_EL RETIRO",       14, "LP1912", "31/12/2025"),
    @(838, "10:42:40", "11:01", "17_ROMERO",          19, "LP1912", "31/12/2025"),
    @(839, "10:42:40", "11:03", "23_HERNANDEZ",       21, "LP1912", "31/12/2025"),
    @(840, "10:42:40", "11:04", "14_ABASTO",          22, "LP1912", "31/12/2025"),
    @(841, "10:42:40", "11:09", "16_SANTA ANA",       27, "LP1912", "31/12/2025"),
    @(842, "10:42:40", "11:11", "15_ABASTO",          29, "LP1912", "31/12/2025"),
    @(843, "10:42:40", "11:25", "16_P MOR-SANTA ANA", 43, "LP1912", "31/12/2025"),
    @(844, "10:42:40", "11:27", "10_OLMOS",           45, "LP1912", "31/12/2025"),
    @(845, "10:42:40", "11:30", "15X38_ABASTO",       48, "LP1912", "31/12/2025"),
    @(846, "10:42:40", "11:33", "23_HERNANDEZ",       51, "LP1912", "31/12/2025"),
    @(847, "10:42:40", "11:34", "10_OLMOS",           52, "LP1912", "31/12/2025"),
    @(848, "10:42:40", "11:40", "215A_EL PATO",       58, "LP1912", "31/12/2025"),
    @(849, "10:42:40", "11:45", "16_SANTA ANA",       63, "LP1912", "31/12/2025"),
    @(850, "10:42:40", "11:53", "15_ABASTO",          71, "LP1912", "31/12/2025"),
    @(851, "10:42:40", "11:54", "225_GOMEZ",          72, "LP1912", "31/12/2025"),
    @(852, "10:42:40", "12:05", "23_HERNANDEZ",       83, "LP1912", "31/12/2025")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": columns A(-) B=Fecha C=Hora_Scrap D=Hora_Llegada
#                      E=Linea F=Minutos G=Parada
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = $updated
$ws2.Range("A3").Value = "Total filas: 62"

$ws2.Cells.Item(63, 2).Value = "31/12/2025"
$ws2.Cells.Item(63, 3).Value = "10:42:40"
$ws2.Cells.Item(63, 4).Value = "11:40"
$ws2.Cells.Item(63, 5).Value = "215A_EL PATO"
$ws2.Cells.Item(63, 6).Value = 58
$ws2.Cells.Item(63, 7).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173": columns A(-) B=Fecha C=Hora_Scrap D=Hora_Llegada
#                     E=Linea F=Minutos G=Parada
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = $updated
$ws3.Range("A3").Value = "Total filas: 101"

$ws3.Cells.Item(102, 2).Value = "31/12/2025"
$ws3.Cells.Item(102, 3).Value = "10:42:45"
$ws3.Cells.Item(102, 4).Value = "11:43"
$ws3.Cells.Item(102, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(102, 6).Value = 61
$ws3.Cells.Item(102, 7).Value = "L6203"
